{"js": "// Replace the three-digit \u00f7 one-digit division answers throughout the\n// document. Each \"before\" string occurs exactly once in the document, so a\n// simple case-sensitive search + whole-text replace is safe and unambiguous.\nconst replacements = [\n  [\"206\u00f77=29, 3\", \"790\u00f78=98, 6\"],\n  [\"206\u00f78=25, 6\", \"142\u00f72=71, 0\"],\n  [\"821\u00f76=136, 5\", \"180\u00f79=20, 0\"],\n  [\"340\u00f77=48, 4\", \"132\u00f72=66, 0\"],\n  [\"394\u00f76=65, 4\", \"332\u00f73=110, 2\"],\n  [\"177\u00f75=35, 2\", \"829\u00f76=138, 1\"],\n  [\"570\u00f79=63, 3\", \"378\u00f72=189, 0\"],\n  [\"265\u00f74=66, 1\", \"516\u00f73=172, 0\"],\n  [\"663\u00f72=331, 1\", \"275\u00f78=34, 3\"],\n  [\"168\u00f77=24, 0\", \"981\u00f77=140, 1\"],\n  [\"775\u00f75=155, 0\", \"225\u00f77=32, 1\"],\n  [\"706\u00f76=117, 4\", \"452\u00f76=75, 2\"],\n  [\"279\u00f76=46, 3\", \"886\u00f79=98, 4\"],\n  [\"519\u00f72=259, 1\", \"845\u00f77=120, 5\"],\n  [\"984\u00f75=196, 4\", \"167\u00f79=18, 5\"],\n  [\"120\u00f72=60, 0\", \"786\u00f76=131, 0\"],\n  [\"635\u00f76=105, 5\", \"642\u00f72=321, 0\"],\n  [\"782\u00f72=391, 0\", \"761\u00f75=152, 1\"],\n  [\"807\u00f76=134, 3\", \"613\u00f72=306, 1\"],\n  [\"929\u00f75=185, 4\", \"963\u00f72=481, 1\"],\n  [\"120\u00f76=20, 0\", \"660\u00f75=132, 0\"],\n  [\"101\u00f78=12, 5\", \"718\u00f78=89, 6\"],\n  [\"948\u00f72=474, 0\", \"687\u00f74=171, 3\"],\n  [\"933\u00f79=103, 6\", \"584\u00f75=116, 4\"],\n  [\"104\u00f76=17, 2\", \"915\u00f79=101, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit \u00f7 one-digit division answers throughout the\n# document. Each \"before\" string occurs exactly once, so a plain\n# Find/Replace (wdReplaceAll) pass per pair is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"206\u00f77=29, 3\", \"790\u00f78=98, 6\"),\n    @(\"206\u00f78=25, 6\", \"142\u00f72=71, 0\"),\n    @(\"821\u00f76=136, 5\", \"180\u00f79=20, 0\"),\n    @(\"340\u00f77=48, 4\", \"132\u00f72=66, 0\"),\n    @(\"394\u00f76=65, 4\", \"332\u00f73=110, 2\"),\n    @(\"177\u00f75=35, 2\", \"829\u00f76=138, 1\"),\n    @(\"570\u00f79=63, 3\", \"378\u00f72=189, 0\"),\n    @(\"265\u00f74=66, 1\", \"516\u00f73=172, 0\"),\n    @(\"663\u00f72=331, 1\", \"275\u00f78=34, 3\"),\n    @(\"168\u00f77=24, 0\", \"981\u00f77=140, 1\"),\n    @(\"775\u00f75=155, 0\", \"225\u00f77=32, 1\"),\n    @(\"706\u00f76=117, 4\", \"452\u00f76=75, 2\"),\n    @(\"279\u00f76=46, 3\", \"886\u00f79=98, 4\"),\n    @(\"519\u00f72=259, 1\", \"845\u00f77=120, 5\"),\n    @(\"984\u00f75=196, 4\", \"167\u00f79=18, 5\"),\n    @(\"120\u00f72=60, 0\", \"786\u00f76=131, 0\"),\n    @(\"635\u00f76=105, 5\", \"642\u00f72=321, 0\"),\n    @(\"782\u00f72=391, 0\", \"761\u00f75=152, 1\"),\n    @(\"807\u00f76=134, 3\", \"613\u00f72=306, 1\"),\n    @(\"929\u00f75=185, 4\", \"963\u00f72=481, 1\"),\n    @(\"120\u00f76=20, 0\", \"660\u00f75=132, 0\"),\n    @(\"101\u00f78=12, 5\", \"718\u00f78=89, 6\"),\n    @(\"948\u00f72=474, 0\", \"687\u00f74=171, 3\"),\n    @(\"933\u00f79=103, 6\", \"584\u00f75=116, 4\"),\n    @(\"104\u00f76=17, 2\", \"915\u00f79=101, 6\")\n)\n\nforeach ($pair in $replacements) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $after\n    $find.Execute($before, $false, $false, $false, $false, $false, $true, 1, $false, $after, 2)\n}\n"}
